# Apply "想去人数" (F column) count updates to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 375
    "F4"  = 32
    "F12" = 1139
    "F17" = 100
    "F24" = 1679
    "F30" = 96
    "F31" = 3977
    "F35" = 1016
    "F39" = 114
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
